$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "SOCKET"

$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 4
$ws.Range("C8").Value = 5
$ws.Range("C9").Value = 5
$ws.Range("C10").Value = 6

$ws.Range("D3").Select()
